$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Tag 2 measurements entered into column G (Anzahl Biegungen) for rows 26-49
$ws.Range("G26").Value = 26
$ws.Range("G27").Value = 9
$ws.Range("G28").Value = 18
$ws.Range("G29").Value = 9
$ws.Range("G30").Value = 8
$ws.Range("G31").Value = 42
$ws.Range("G32").Value = 29
$ws.Range("G33").Value = 31
$ws.Range("G34").Value = 9
$ws.Range("G35").Value = 13
$ws.Range("G36").Value = 15
$ws.Range("G37").Value = 9
$ws.Range("G38").Value = 10
$ws.Range("G39").Value = 9
$ws.Range("G40").Value = 29
$ws.Range("G41").Value = 9
$ws.Range("G42").Value = 20
$ws.Range("G43").Value = 16
$ws.Range("G44").Value = 10
$ws.Range("G45").Value = 8
$ws.Range("G46").Value = 14
$ws.Range("G47").Value = 18
$ws.Range("G48").Value = 12
$ws.Range("G49").Value = 8

# Mark the end of the "Tag 1" block (row 13) and the Block-3 mid-point (row 37)
# with a thin bottom border, like the other block-separator rows in the sheet.
$ws.Range("A13:H13").Borders.Item(9).LineStyle = 1
$ws.Range("A13:H13").Borders.Item(9).Weight = 2

$ws.Range("A37:H37").Borders.Item(9).LineStyle = 1
$ws.Range("A37:H37").Borders.Item(9).Weight = 2

# Move the selection / scroll position to reflect where entry left off
$ws.Range("G50").Select()
